# Applies the "Updated cryptos list" data refresh:
# updates Price (D) / Volume(1h) (E) figures for most rows, and
# swaps the Decentraland / EnergySwap rows (45 <-> 46) including their
# name, link, price and volume columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, exactly as stored in the source feed,
# without letting Excel auto-convert number-looking strings (e.g. "0.5078",
# "1.870.84") into numeric/date values. NumberFormat is restored afterwards
# so cell formatting stays identical to the original workbook.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$updates = [ordered]@{
    "D2" = "30.000.39"
    "E2" = "  -0.23%  "
    "D3" = "1.870.84"
    "E3" = "  -2.69%  "
    "E4" = "  +0.04%  "
    "D5" = "319.22"
    "E5" = "  -2.31%  "
    "E6" = "  +0.04%  "
    "D7" = "0.5078"
    "E7" = "  -3.08%  "
    "D8" = "0.3941"
    "E8" = "  -3.07%  "
    "D9" = "0.08180"
    "E9" = "  -3.47%  "
    "D10" = "42.16"
    "E10" = "  -1.74%  "
    "D11" = "1.094"
    "E11" = "  -3.19%  "
    "D12" = "22.74"
    "E12" = "  +2.44%  "
    "D13" = "1.875.03"
    "E13" = "  -2.57%  "
    "E14" = "  -1.76%  "
    "D15" = "7.176"
    "E15" = "  -2.84%  "
    "D16" = "1.002"
    "E16" = "  +0.06%  "
    "D17" = "92.04"
    "E17" = "  -4.40%  "
    "D18" = "0.00001086"
    "E18" = "  -2.75%  "
    "D19" = "0.06431"
    "E19" = "  -4.61%  "
    "E20" = "  -1.97%  "
    "E21" = "  +0.06%  "
    "D22" = "29.983.69"
    "E22" = "  -0.34%  "
    "D23" = "5.808"
    "E23" = "  -4.27%  "
    "D24" = "11.12"
    "E24" = "  -1.41%  "
    "D25" = "2.147"
    "E25" = "  -2.34%  "
    "D26" = "2.083.53"
    "E26" = "  -2.90%  "
    "D27" = "160.92"
    "E27" = "  +0.19%  "
    "D28" = "20.95"
    "E28" = "  -1.28%  "
    "D29" = "2.236"
    "E29" = "  -9.27%  "
    "D30" = "127.07"
    "E30" = "  -1.40%  "
    "D31" = "1.063"
    "E31" = "  -2.09%  "
    "D32" = "0.1036"
    "E32" = "  -2.43%  "
    "D33" = "5.904"
    "E33" = "  -3.38%  "
    "D34" = "3.729"
    "E34" = "  +1.84%  "
    "D35" = "0.02422"
    "E35" = "  -4.10%  "
    "D36" = "5.224"
    "E36" = "  +0.29%  "
    "D37" = "0.06347"
    "D39" = "1.173"
    "E39" = "  -5.36%  "
    "D40" = "8.487"
    "E40" = "  -6.26%  "
    "D41" = "0.6309"
    "E41" = "  -4.23%  "
    "D42" = "11.23"
    "E42" = "  -3.59%  "
    "D43" = "1.193"
    "E43" = "  -4.39%  "
    "E44" = "  +0.06%  "
    "B45" = "Decentraland"
    "C45" = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
    "D45" = "0.5900"
    "E45" = "  -4.70%  "
    "B46" = "EnergySwap"
    "C46" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D46" = "12.90"
    "E46" = "  -2.29%  "
    "D47" = "3.630"
    "E47" = "  -3.45%  "
    "D48" = "1.995"
    "E48" = "  -4.04%  "
    "D49" = "122.66"
    "E49" = "  -2.62%  "
    "D50" = "1.198"
    "E50" = "  -3.77%  "
    "D51" = "1.125"
    "E51" = "  -3.44%  "
}

foreach ($ref in $updates.Keys) {
    Set-TextValue $ref $updates[$ref]
}

